$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("15")
$ws.Range("E2").Value = 8509.227
$ws.Range("I2").Value = -0.0191543497943117
$ws.Range("J2").Value = 0.0255241288994157
$ws.Range("K2").Value = 0.0191543497943117
$ws.Range("L2").Value = -1.59304907081334
$ws.Range("F3").Value = 0.274027142180158
$ws.Range("H3").Value = 52.302625
$ws.Range("I3").Value = -0.377725427512456
$ws.Range("J3").Value = 0.517554434954518
$ws.Range("K3").Value = 0.377725427512456
$ws.Range("L3").Value = -0.286043965567491
$ws.Range("E7").Value = 8757.557
$ws.Range("I7").Value = -0.141260996271992
$ws.Range("J7").Value = 0.173934288769413
$ws.Range("K7").Value = 0.141260996271992
$ws.Range("L7").Value = -0.759614794331786
$ws.Range("E8").Value = 4026.08
$ws.Range("I8").Value = 1.40185787623314
$ws.Range("J8").Value = 1.2518618359821
$ws.Range("K8").Value = 1.40185787623314
$ws.Range("L8").Value = 0.0975563998155705

$ws = $wb.Worksheets.Item("3")
$ws.Range("F3").Value = 0.514403686709471
$ws.Range("H3").Value = 52.1867083333333
$ws.Range("I3").Value = -0.0148789334135516
$ws.Range("J3").Value = 0.0204362129975192
$ws.Range("K3").Value = 0.0148789334135516
$ws.Range("L3").Value = -1.6895995795088

$ws = $wb.Worksheets.Item("5")
$ws.Range("E3").Value = 3799.521
$ws.Range("F3").Value = 0.216853459021603
$ws.Range("H3").Value = 53.521625
$ws.Range("I3").Value = -0.1288164418135
$ws.Range("J3").Value = 0.17242135645319
$ws.Range("K3").Value = 0.1288164418135
$ws.Range("L3").Value = -0.76340894259866

$ws = $wb.Worksheets.Item("5a")
$ws.Range("F3").Value = 0.519697920654582
$ws.Range("H3").Value = 53.17225
$ws.Range("I3").Value = -0.197842978098901
$ws.Range("J3").Value = 0.266638686605334
$ws.Range("K3").Value = 0.197842978098901
$ws.Range("L3").Value = -0.574076838566927

$ws = $wb.Worksheets.Item("6")
$ws.Range("E2").Value = 12349.884
$ws.Range("I2").Value = -0.0248130778879441
$ws.Range("J2").Value = 0.0290951118423101
$ws.Range("K2").Value = 0.0248130778879441
$ws.Range("L2").Value = -1.53617996902892
$ws.Range("F4").Value = 0.391295963011507
$ws.Range("H4").Value = 53.2519166666667
$ws.Range("I4").Value = -0.0143649998659057
$ws.Range("J4").Value = 0.0193268918419151
$ws.Range("K4").Value = 0.0143649998659057
$ws.Range("L4").Value = -1.71383798375942

$ws = $wb.Worksheets.Item("6a")
$ws.Range("F2").Value = 0.303866762065075
$ws.Range("A4").Value = 1.2025711866528
$ws.Range("F4").Value = 0.23937442308451
$ws.Range("G4").Value = 15.943041869971
$ws.Range("H4").Value = 52.15875
$ws.Range("I4").Value = -0.330011044031721
$ws.Range("J4").Value = 0.453426571314305
$ws.Range("K4").Value = 0.330011044031721
$ws.Range("L4").Value = -0.343493033311727
$ws.Range("E5").Value = 1573.128
$ws.Range("I5").Value = -0.493194293879235
$ws.Range("J5").Value = 0.440035837567019
$ws.Range("K5").Value = 0.493194293879235
$ws.Range("L5").Value = -0.356511952096092

$ws = $wb.Worksheets.Item("7")
$ws.Range("A3").Value = 1.92653479913752
$ws.Range("E3").Value = 7475.905
$ws.Range("F3").Value = 0.327860796281033
$ws.Range("G3").Value = 84.4373897169224
$ws.Range("H3").Value = 52.7094583333333
$ws.Range("I3").Value = -0.552638322840444
$ws.Range("J3").Value = 0.751239608585646
$ws.Range("K3").Value = 0.552638322840444
$ws.Range("L3").Value = -0.124221522263913

$ws = $wb.Worksheets.Item("9")
$ws.Range("E2").Value = 3903.276
$ws.Range("I2").Value = -0.0265712094335327
$ws.Range("J2").Value = 0.037428936539669
$ws.Range("K2").Value = 0.0265712094335327
$ws.Range("L2").Value = -1.42679251222607
$ws.Range("F4").Value = 0.196666010534818
$ws.Range("H4").Value = 50.6177916666667
$ws.Range("I4").Value = -0.145662173571032
$ws.Range("J4").Value = 0.206182707231445
$ws.Range("K4").Value = 0.145662173571032
$ws.Range("L4").Value = -0.685747762277013
